$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I/J: advective-flux-at-zinf variants of the existing C/D calcs ---

# Row 18: loc_fPOC-like term divided differently -> 1/(1-por)*... (mirrors C18/D18 but inverted porosity factor)
$i18 = $ws.Range("I18")
$i18.Formula = "=1/(1-C11)*C15/(C7)"
$i18.Font.Name = "Arial"
$i18.Font.Size = 10
$i18.Font.Bold = $false

$j18 = $ws.Range("J18")
$j18.Formula = "=1/(1-D11)*D15/(D7)"
$j18.Font.Name = "Arial"
$j18.Font.Size = 10
$j18.Font.Bold = $false

# Row 20: new maroon bold header note above the wt% block
$h20 = $ws.Range("H20")
$h20.Value = "HERE IS WHAT I, WHEN DIVIDE WITH (1-por):"
$h20.Font.Name = "Arial"
$h20.Font.Size = 10
$h20.Font.Bold = $true
$h20.Font.Color = 128

# Row 21: wt% using density dum_den, for the new I/J columns (same look as C21/D21)
$i21 = $ws.Range("I21")
$i21.Formula = "=100*I18*12/C12"
$i21.Font.Bold = $true
$i21.Font.Color = 255
$i21.NumberFormat = "0.00"

$j21 = $ws.Range("J21")
$j21.Formula = "=100*J18*12/D12"
$j21.Font.Bold = $true
$j21.Font.Color = 255
$j21.NumberFormat = "0.00"

# Row 22: wt% using density of 1 (as in GENIE), for the new I/J columns (same look as C22/D22)
$i22 = $ws.Range("I22")
$i22.Formula = "=100*I18*12/1"
$i22.Font.Bold = $true
$i22.NumberFormat = "0.00"

$j22 = $ws.Range("J22")
$j22.Formula = "=100*J18*12/1"
$j22.Font.Bold = $true
$j22.NumberFormat = "0.00"

# Column H gets a bit wider to fit the new header text
$ws.Columns.Item(8).ColumnWidth = 14.1

# Restore the selection to where the author left off editing
$ws.Range("D29").Select()
